$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells
$ws.Range("E2").Value = "white"
$ws.Range("B3").Value = "Diego"
$ws.Range("E3").Value = "red"
$ws.Range("E4").Value = "white"

# Add new rows 5 and 6
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "awdas"
$ws.Range("C5").Value = "aaa"
$ws.Range("D5").Value = "'123"
$ws.Range("E5").Value = "white"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "giulio"
$ws.Range("C6").Value = "abacate"
$ws.Range("D6").Value = "123aa"
$ws.Range("E6").Value = "white"
